$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Merge the three runs that make up the tail of the "Scenario:"
#    paragraph ("<space>", "You have", "<space>set an ... shows-up.")
#    into a single run, without changing the visible text.
# ------------------------------------------------------------------

# Locate the end of the literal "Scenario:" label - the merge target
# starts right after it and runs to the end of the paragraph (but
# before the paragraph mark).
$findRng = $d.Content
$null = $findRng.Find.Execute("Scenario:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeStart = $findRng.End

$paraRng = $d.Range($mergeStart, $mergeStart)
$paraRng.Expand(4) | Out-Null   # wdParagraph -> grow to the enclosing paragraph
$mergeEnd = $paraRng.End - 1    # exclude the trailing paragraph mark

$target = $d.Range($mergeStart, $mergeEnd)

$leftQuote = [char]0x2018
$rightQuote = [char]0x2019
$mergedText = " You have set an Out of Office notification in your e-mail box. In the Out of Office message, you will mention that if the message is urgent, people can send an e-mail containing just one exclamation point " + $leftQuote + "!" + $rightQuote + ". You also want to get a mobile notification when such message shows-up."

# Minimal XML-escape in case the literal text above is ever changed.
$mergedTextXml = $mergedText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="0019443B"><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Segoe UI" w:hAnsi="Segoe UI" w:cs="Times New Roman"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">$mergedTextXml</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$target.InsertXML($xml)

# ------------------------------------------------------------------
# 2. Mark the existing comment (paraId 76A818EB) as "Done".
# ------------------------------------------------------------------
$comments = $d.Comments
for ($i = 1; $i -le $comments.Count; $i++) {
    $comments.Item($i).Done = $true
}
